$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.814.16"
$ws.Range("E2").Value = "  +8.44%  "
$ws.Range("D3").Value = "3.218.02"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'397.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "'109.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.57%  "
$ws.Range("E7").Value = "  +3.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("D10").Value = "'39.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").Value = "'0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.23%  "
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").Value = "3.725.04"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "'19.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "'8.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.65%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.222.89"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'1.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.39%  "
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "55.693.76"
$ws.Range("E19").Value = "  +8.10%  "
$ws.Range("E20").Value = "  +3.83%  "
$ws.Range("D21").Value = "'0.0000103"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.64%  "
$ws.Range("E22").Value = "  +5.95%  "
$ws.Range("D23").Value = "'301.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.16%  "
$ws.Range("D24").Value = "'75.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.70%  "
$ws.Range("D25").Value = "'3.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "'28.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.95%  "
$ws.Range("D28").Value = "'7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.16%  "
$ws.Range("D29").Value = "'0.175"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("E32").Value = "  +9.33%  "
$ws.Range("D33").Value = "'0.0492"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").Value = "'36.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "'51.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("E37").Value = "  +23.87%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").Value = "'134.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("D41").Value = "'4.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.80%  "
$ws.Range("D42").Value = "'1.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").Value = "'17.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.119"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "'22.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'2.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +45.24%  "
$ws.Range("D49").Value = "'2.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "2.138.04"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "'0.0363"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.07%  "
